# Replace the trailing sentence of the PL/SQL-functions paragraph with a new
# paragraph/sentence, leaving the single leading space from the old sentence
# in place, and stamp a "_GoBack" bookmark (Word's "last edit location"
# marker) right after the newly typed text.
#
#   ... are required. In user defined functions if not anonymous can be
#   called anytime in SQL queries.
#
# becomes
#
#   ... are required.
#   We will soon discuss with examples on user defined functions in oracle.

$d = $word.ActiveDocument

$oldSentence = "In user defined functions if not anonymous can be called anytime in SQL queries."
$newSentence = "We will soon discuss with examples on user defined functions in oracle."

# --- Locate the sentence to remove -----------------------------------------
$target = $d.Content
$found = $target.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetStart = $target.Start

# Figure out which paragraph holds the match, so fresh Paragraph objects can
# be re-fetched by index after each edit below (keeping everything in sync).
$paraIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $targetStart -lt $p.Range.End) {
        $paraIndex = $i
    }
}

# --- Delete the sentence, keeping the preceding run (the lone space) intact-
# The character right before the match is the leading space that must
# survive as its own run. Toggling its formatting on/off around the deletion
# stops the engine from silently re-merging that now-lone space into the
# previous (identically formatted) run.
$spaceRange = $d.Range($targetStart - 1, $targetStart)
$spaceRange.Bold = 1
$target.Delete()
$spaceRange = $d.Range($targetStart - 1, $targetStart)
$spaceRange.Bold = 0

# --- Split a new paragraph off right after it -------------------------------
$para = $d.Paragraphs.Item($paraIndex)
$splitPoint = $d.Range($para.Range.End - 1, $para.Range.End - 1)
$splitPoint.InsertParagraphAfter()

# --- Type the replacement sentence into the new (empty) paragraph ----------
# A throw-away trailing marker character is appended for now; see below.
$newPara = $d.Paragraphs.Item($paraIndex + 1)
$newRange = $newPara.Range
$newRange.InsertAfter($newSentence + "#")

# --- Stamp the "_GoBack" bookmark Word leaves at the last edit location ----
# A zero-length bookmark placed at the very end of a paragraph's text (i.e.
# immediately before the paragraph mark, with nothing following it) would
# render *before* the preceding run instead of after it, so the throw-away
# "#" above gives it something to split against; once the bookmark is
# anchored, that placeholder character is deleted again.
$newPara = $d.Paragraphs.Item($paraIndex + 1)
$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
